$d = $word.ActiveDocument

# --- Edit 1: the empty paragraph right after "What do you mean by BI? Explain."
# had a stray _GoBack bookmark (left behind from a prior save); drop it while
# keeping the paragraph itself (and its formatting) untouched.
$target1 = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.Trim() -eq "") {
        $prevText = $d.Paragraphs($i - 1).Range.Text
        if ($prevText -like "*What do you mean by BI*") {
            $target1 = $para
            break
        }
    }
}
if ($target1 -ne $null) {
    $xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:highlight w:val="yellow"/></w:rPr></w:pPr></w:p>'
    $target1.Range.InsertXML($xml1)
}

# --- Edit 2: the _GoBack bookmark now belongs at the end of the document's
# last real content paragraph (the one ending "...built-in algorithms."), and
# a new graded Q&A section (question 5 + answers) follows it.
$target2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*built-in algorithms.*") {
        $target2 = $para
        break
    }
}
if ($target2 -ne $null) {
    $xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00286434" w:rsidRDefault="00286434" w:rsidP="00286434"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="420" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="181818"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="181818"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>At your company, you can use prescriptive analytics to conduct manual analyses, develop proprietary algorithms, or use third-party analytics tools with built-in algorithms.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="181818"/><w:sz w:val="36"/><w:szCs w:val="26"/><w:highlight w:val="yellow"/></w:rPr><w:t>5.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Write five real-life questions that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/><w:highlight w:val="yellow"/></w:rPr><w:t>PowerBi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> can solve.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t>1. Sales Performance Analysis: How are sales performing across different regions, products, or time periods? Which sales channels are generating the highest revenue? Power BI can help analyze sales data to identify trends, top-selling products, and opportunities for improvement.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">2. Customer Segmentation: How can customers be segmented based on demographics, buying behavior, or preferences? Power BI can integrate customer data from multiple sources and provide insights on </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:lastRenderedPageBreak/><w:t>customer segments to tailor marketing strategies and improve customer satisfaction.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t>3. Supply Chain Optimization: How can the supply chain be optimized to reduce costs and improve efficiency? Power BI can analyze inventory levels, transportation data, and production metrics to identify bottlenecks, optimize inventory management, and streamline supply chain operations.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t>4. Financial Performance Monitoring: How is the company''s financial performance? Power BI can consolidate financial data, generate key financial ratios, track revenue and expenses, and provide real-time insights into profitability, cash flow, and financial health.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t>5. Social Media Analytics: How is the company''s brand being perceived on social media? Power BI can integrate social media data and sentiment analysis to track brand mentions, sentiment trends, customer feedback, and measure the effectiveness of social media campaigns.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:eastAsia="Times New Roman" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:sz w:val="28"/><w:szCs w:val="21"/></w:rPr><w:t>Please note that Power BI is a versatile tool that can solve a wide range of analytical challenges across various industries. These are just a few examples, and the possibilities are extensive based on the specific data and business needs.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="420" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="181818"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr></w:p>'
    $target2.Range.InsertXML($xml2)
}
